$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A94").Value = "2025/12/06 17:00"
$ws.Range("B94").Value = "-"
$ws.Range("C94").Value = "-"
$ws.Range("D94").Value = "-"
$ws.Range("E94").Value = "-"
$ws.Range("F94").Value = "-"
$ws.Range("G94").Value = "-"
